# Weekly update: insert two new daily price records for Ají (Americana (o))
# at Femacal de La Calera, pushing the existing rows 367-391 down to 369-393.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of this block (row 367), shifting the
# previously-existing rows 367:391 down to 369:393.
$ws.Rows("367:368").Insert()

# New row 367
$ws.Range("A367").Value = 3
$ws.Range("B367").Value = "Femacal de La Calera"
$ws.Range("C367").Value = "Coquimbo"
$ws.Range("D367").Value = 44610
$ws.Range("E367").Value = 5
$ws.Range("F367").Value = 100112021
$ws.Range("G367").Value = "Ají"
$ws.Range("H367").Value = "Americana (o)"
$ws.Range("I367").Value = "Primera"
$ws.Range("J367").Value = 73
$ws.Range("K367").Value = 14000
$ws.Range("L367").Value = 15000
$ws.Range("M367").Value = 14521
$ws.Range("N367").Value = "$/caja 15 kilos"
$ws.Range("O367").Value = "Limache"
$ws.Range("P367").Value = 968
$ws.Range("Q367").Value = 15
$ws.Range("R367").Value = "Hortaliza"

# New row 368
$ws.Range("A368").Value = 3
$ws.Range("B368").Value = "Femacal de La Calera"
$ws.Range("C368").Value = "Coquimbo"
$ws.Range("D368").Value = 44610
$ws.Range("E368").Value = 5
$ws.Range("F368").Value = 100112021
$ws.Range("G368").Value = "Ají"
$ws.Range("H368").Value = "Americana (o)"
$ws.Range("I368").Value = "Primera"
$ws.Range("J368").Value = 45
$ws.Range("K368").Value = 24000
$ws.Range("L368").Value = 25000
$ws.Range("M368").Value = 24444
$ws.Range("N368").Value = "$/caja 25 kilos"
$ws.Range("O368").Value = "Provincia de Limarí"
$ws.Range("P368").Value = 978
$ws.Range("Q368").Value = 25
$ws.Range("R368").Value = "Hortaliza"
